# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 26;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 37;  I = "ba"; J = "Appreciation" },
    @{ Row = 39;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 42;  I = "ba"; J = "Appreciation" },
    @{ Row = 47;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 50;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 89;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 103; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 112; I = "sv"; J = "Statement-opinion" },
    @{ Row = 119; I = "sv"; J = "Statement-opinion" },
    @{ Row = 121; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 126; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 132; I = "sv"; J = "Statement-opinion" },
    @{ Row = 136; I = "sv"; J = "Statement-opinion" },
    @{ Row = 140; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 148; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 155; I = "%"; J = "Uninterpretable" },
    @{ Row = 162; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 169; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.I
    $ws.Range("J" + $u.Row).Value = $u.J
}
